$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-65 (Cluster Name, Active Cases). Row 1 (header) is untouched.
$rows = @(
    @(2, '3364 Assisi Centre Aged Care Rosanna', 29),
    @(3, '3376 Royal Freemasons Coppin Centre Melbourne', 21),
    @(4, '3622 Olivet Care Aged Care Services Ringwood', 12),
    @(5, '3825 TLC Forest Lodge Residential Aged Care Frankston North', 14),
    @(6, '3961 Heritage Care Water Gardens Aged Care Facility Sydenham', 22),
    @(7, '4167 Royal Freemasons Centennial Lodge Wantirna South', 23),
    @(8, '4282 Villa Maria Catholic Homes (VMCH) Wantirna At-Home Aged Care Stud Road Wantirna', 11),
    @(9, 'AG Industries Pty Ltd Factory Thomastown', 15),
    @(10, 'Alamanda K9 College Point Cook', 11),
    @(11, 'Antonine College Cedar Campus Coburg', 12),
    @(12, 'Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh', 33),
    @(13, 'Baden Powell College Tarneit', 12),
    @(14, 'Covenant College Bell Post Hill', 26),
    @(15, 'Dandenong South Primary School Dandenong', 11),
    @(16, 'Devon Meadows Primary School Devon Meadows', 12),
    @(17, 'Flemington Racecourse Flemington', 12),
    @(18, 'Gilly''s Early Learning Centre Balaclava', 11),
    @(19, 'Gladstone Park Secondary College 29 Oct Gladstone Park', 10),
    @(20, 'Guardian Childcare & Education Moorabbin', 13),
    @(21, 'Hamlyn Banks Primary School Hamlyn Heights', 10),
    @(22, 'Hazel Glen College Doreen', 14),
    @(23, 'Hazelwood North Primary School Hazelwood North', 29),
    @(24, 'Islamic College of Melbourne Tarneit Oct Nov', 30),
    @(25, 'Lyndhurst Primary School Lyndhurst', 13),
    @(26, 'Master Poultry Group West Footscray', 14),
    @(27, 'Minaret College Officer Campus Officer', 17),
    @(28, 'Morwell Park Primary School Morwell Outbreak', 58),
    @(29, 'Narre Warren South P-12 College Narre Warren South', 13),
    @(30, 'Nido Early School Woodend', 12),
    @(31, 'Northern Bay College Wexford Campus Corio', 25),
    @(32, 'Northern Health Northern Hospital Epping Emergency Department Tier 1B', 14),
    @(33, 'Northern Health The Northern Hospital Epping', 15),
    @(34, 'Oakleigh South Primary School Oakleigh South', 13),
    @(35, 'Pentland Primary School Darley', 11),
    @(36, 'Rutherglen Motor Inn and Walkabout Motel Rutherglen', 22),
    @(37, 'Sirius College Ibrahim Dellal Campus Sunshine', 11),
    @(38, 'Sirius College Shepparton Campus Shepparton', 22),
    @(39, 'Smartie Pants Early Learning and Development Diamond Creek', 21),
    @(40, 'Social Gathering Woodvale 30 Oct', 10),
    @(41, 'Society Restaurant Melbourne', 26),
    @(42, 'St Ambrose Parish Primary School Woodend', 11),
    @(43, 'St Brendans Primary School Shepparton', 10),
    @(44, 'St Clare''s Primary School Officer', 12),
    @(45, 'St Georges Road Primary School Shepparton', 15),
    @(46, 'St Joseph''s School Quarry Hill', 32),
    @(47, 'St Louis de Montfort''s School Aspendale', 13),
    @(48, 'St Paul''s Primary School Sunshine West', 15),
    @(49, 'St Vincents Hospital Melbourne Emergency Department Fitzroy', 14),
    @(50, 'Stevensville Primary School St Albans', 10),
    @(51, 'Stockdale Road Primary School Traralgon', 33),
    @(52, 'Supreme Caravans Manufacturing Campbellfield', 19),
    @(53, 'Templestowe Park Primary School Templestowe', 31),
    @(54, 'The Lake Primary School Cabarita', 24),
    @(55, 'The Royal Children''s Hospital Melbourne Emergency Department Parkville', 15),
    @(56, 'Top Yard Rooftop Melbourne', 13),
    @(57, 'Truganina P-9 College Truganina', 10),
    @(58, 'Tucker Road Bentleigh Primary School Bentleigh', 11),
    @(59, 'Warragul Regional College Warragul', 19),
    @(60, 'Werribee Mercy Hospital Emergency Department', 18),
    @(61, 'Western Health Sunshine Hospital Emergency Department St Albans', 13),
    @(62, 'Wodonga Primary School Wodonga', 13),
    @(63, 'Wodonga South Primary School Wodonga', 10),
    @(64, 'Wyndham Christian College Wyndham Vale', 11),
    @(65, 'Yeshivah College St Kilda East', 23)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $name = $r[1]
    $val = $r[2]
    $ws.Cells.Item($rowNum, 1).Value = $name
    $ws.Cells.Item($rowNum, 2).Value = $val
}

# Rows 66-69 no longer exist in the updated data - remove them.
$ws.Rows("66:69").Delete()
